# Arrumado a escrita da apresentação
$p = $ppt.ActivePresentation

# Slide 3: "Como toda criança, meu interesse começou na infância."
# -> split into two runs: "Como todo mundo, " + "meu interesse começou na infância."
$s3 = $p.Slides.Item(3)
$shape3 = $s3.Shapes.Item(4)
$tr3 = $shape3.TextFrame.TextRange
$prefix3 = $tr3.Characters(1, 19)
$prefix3.Text = "Como todo mundo, "

# Slide 4: "Por meio de minhas ações, sendo como motivador as lições ensinadas,
#           no qual ajudam no desenvolvimento pessoal."
# -> "Por meio de minhas ações, sendo como motivador, as lições ensinadas,
#     no qual ajudaram no meu desenvolvimento pessoal."
$s4 = $p.Slides.Item(4)
$shape4 = $s4.Shapes.Item(4)
$tr4 = $shape4.TextFrame.TextRange
$tr4.Text = "Por meio de minhas ações, sendo como motivador, as lições ensinadas, no qual ajudaram no meu desenvolvimento pessoal."
